$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (rich-text shared strings): bulletin number + week dates
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  18"
$ws.Range("C9").Value = "Report Covering the Week  5/1/2023  Through  5/7/2023"

# ---------------------------------------------------------------------------
# Row 15 (Rape) - Week to Date 2023 count
# ---------------------------------------------------------------------------
$ws.Range("F15").Value = 1

# ---------------------------------------------------------------------------
# Row 16 (Robbery) - reshuffled: a new "Week to Date 2023" column of data
# appears (C16), bumping the previous C/D/E content over to D/E (as N/A
# placeholders), plus refreshed 28-day/YTD/2-year figures.
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 1
$ws.Range("F16").Copy()
$ws.Range("C16").PasteSpecial(-4122)  # xlPasteFormats -> numeric style (15)

$ws.Range("D16").Value = "'0"
$ws.Range("A16").Copy()
$ws.Range("D16").PasteSpecial(-4122)  # xlPasteFormats -> text style (14)

$ws.Range("E16").Value = "***.*"
$ws.Range("A16").Copy()
$ws.Range("E16").PasteSpecial(-4122)  # xlPasteFormats -> text style (14)

$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -77.777777777777
$ws.Range("I16").Value = 22
$ws.Range("K16").Value = -33.333333333333
$ws.Range("L16").Value = 340
$ws.Range("M16").Value = -15.384615384615
$ws.Range("N16").Value = -83.076923076923

# ---------------------------------------------------------------------------
# Row 17 (Fel. Assault)
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 500
$ws.Range("F17").Value = 9
$ws.Range("H17").Value = 80
$ws.Range("I17").Value = 35
$ws.Range("J17").Value = 25
$ws.Range("K17").Value = 40
$ws.Range("L17").Value = 118.75
$ws.Range("M17").Value = 66.666666666666
$ws.Range("N17").Value = -22.222222222222

# ---------------------------------------------------------------------------
# Row 18 (Burglary)
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -16.666666666666
$ws.Range("F18").Value = 22
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = 57.142857142857
$ws.Range("I18").Value = 119
$ws.Range("J18").Value = 106
$ws.Range("K18").Value = 12.264150943396
$ws.Range("L18").Value = 48.75
$ws.Range("M18").Value = 32.222222222222
$ws.Range("N18").Value = -65.706051873198

# ---------------------------------------------------------------------------
# Row 19 (Gr. Larceny)
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 30.769230769230
$ws.Range("F19").Value = 52
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = 23.809523809523
$ws.Range("I19").Value = 258
$ws.Range("J19").Value = 190
$ws.Range("K19").Value = 35.789473684210
$ws.Range("L19").Value = 116.806722689076
$ws.Range("M19").Value = 85.611510791366
$ws.Range("N19").Value = 51.764705882352

# ---------------------------------------------------------------------------
# Row 20 (G.L.A.)
# ---------------------------------------------------------------------------
$ws.Range("F20").Value = 12
$ws.Range("H20").Value = 140
$ws.Range("I20").Value = 52
$ws.Range("J20").Value = 35
$ws.Range("K20").Value = 48.571428571428
$ws.Range("L20").Value = 160
$ws.Range("M20").Value = 10.638297872340
$ws.Range("N20").Value = -95.373665480427

# ---------------------------------------------------------------------------
# Row 21 (TOTAL)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 32
$ws.Range("E21").Value = 45.454545454545
$ws.Range("F21").Value = 98
$ws.Range("G21").Value = 75
$ws.Range("H21").Value = 30.666666666666
$ws.Range("I21").Value = 490
$ws.Range("J21").Value = 389
$ws.Range("K21").Value = 25.964010282776
$ws.Range("L21").Value = 104.166666666667
$ws.Range("M21").Value = 50.769230769230
$ws.Range("N21").Value = -73.106476399560

# ---------------------------------------------------------------------------
# Row 24 (Petit Larceny)
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = -51.724137931034
$ws.Range("F24").Value = 48
$ws.Range("G24").Value = 73
$ws.Range("H24").Value = -34.246575342465
$ws.Range("I24").Value = 189
$ws.Range("J24").Value = 290
$ws.Range("K24").Value = -34.827586206896
$ws.Range("L24").Value = 6.779661016949
$ws.Range("M24").Value = 8.620689655172

# ---------------------------------------------------------------------------
# Row 25 (Misd. Assault)
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 40
$ws.Range("F25").Value = 18
$ws.Range("G25").Value = 11
$ws.Range("H25").Value = 63.636363636363
$ws.Range("I25").Value = 75
$ws.Range("J25").Value = 71
$ws.Range("K25").Value = 5.633802816901
$ws.Range("L25").Value = 97.368421052631
$ws.Range("M25").Value = 44.230769230769

# ---------------------------------------------------------------------------
# Row 26 (UCR Rape*)
# ---------------------------------------------------------------------------
$ws.Range("F26").Value = 1

# ---------------------------------------------------------------------------
# Row 27 (Other Sex Crimes)
# ---------------------------------------------------------------------------
$ws.Range("F27").Value = 1
$ws.Range("H27").Value = 0

# ---------------------------------------------------------------------------
# Row 30 (Hate Crimes)
# ---------------------------------------------------------------------------
$ws.Range("C30").Value = 1
$ws.Range("J30").Copy()
$ws.Range("C30").PasteSpecial(-4122)  # xlPasteFormats -> numeric style (15)

$ws.Range("F30").Value = 2
$ws.Range("I30").Value = 2
$ws.Range("K30").Value = 100
$ws.Range("L30").Value = 0
